$p = $ppt.ActivePresentation

# --- Update the auto "datetimeFigureOut" date fields on the handout
# and notes masters (6/24/2022 -> 4/17/2023) ---
$handoutDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$handoutDate.Text = "4/17/2023"

$notesDate = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDate.Text = "4/17/2023"

# --- Merge the "upervisor" / ": " runs on slide 1 into a single run
# "upervisor: " (same formatting, just collapsing the run split) ---
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange
$fullText = $textRange.Text
$startIdx = $fullText.IndexOf("upervisor") + 1
$subRange = $textRange.Characters($startIdx, 11)
$subRange.Text = "upervisor: "
